$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B10 value to 3
$ws.Range("B10").Value = 3

# Set D10 to the new text (will be stored as shared string)
$ws.Range("D10").Value = "Informa que a partilha foi realizada com sucesso"

# Update column D width
$ws.Columns("D").ColumnWidth = 44.140625

# Update selection to D10
$ws.Range("D10").Select()
